# ReactZero.docx - "Completed lesion 10"
# 1) Drop the stray _GoBack bookmark from the "...render la trang thai quan
#    trong nhat." paragraph (it will be re-added at the very end of the doc).
# 2) Turn the trailing empty paragraph into the new "6. Props" section, ending
#    with a paragraph holding the relocated _GoBack bookmark.

$d = $word.ActiveDocument

# --- Step 1: locate & rewrite the paragraph that owns the "_GoBack" bookmark ---
$bookmarkParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D43556" w:rsidRPr="001479B2" w:rsidRDefault="00D43556" w:rsidP="001479B2"><w:r><w:t>Component có rất nhiều trạng thái nhưng render là trạng thái quan trọng nhất.</w:t></w:r></w:p>
'@

$bookmarkPara = $null
$findRange = $d.Content
$found = $findRange.Find.Execute("quan trọng nhất.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bookmarkPara = $findRange.Paragraphs(1)
} else {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.WordOpenXML -like "*_GoBack*") {
            $bookmarkPara = $p
        }
    }
}
if ($bookmarkPara -ne $null) {
    $null = $bookmarkPara.Range.InsertXML($bookmarkParaXml)
}

# --- Step 2: replace the trailing empty paragraph with the new "6. Props" section ---
$newSectionXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>6. Props</w:t></w:r></w:p><w:p><w:r><w:t>Props là thuộc tính của component. Props không bao giờ thay đổi.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>&lt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:color w:val="000080"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">ThuatNguyen </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:color w:val="0000FF"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:color w:val="008000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t xml:space="preserve">="ReactJS" </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:color w:val="0000FF"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>teacher</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:color w:val="008000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>="Mr.Khoa"</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>React Subject</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>&lt;/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:color w:val="000080"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>ThuatNguyen</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="14"/><w:szCs w:val="14"/><w:shd w:val="clear" w:color="auto" w:fill="EFEFEF"/></w:rPr><w:t>&gt;</w:t></w:r></w:p><w:p><w:r><w:t>Có 1 props đặc biệt là children</w:t></w:r></w:p><w:p><w:r><w:t>Ở đây React Subject là children</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$lastPara = $d.Paragraphs.Last
$null = $lastPara.Range.InsertXML($newSectionXml)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
